$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 6).Value = 169
$ws.Cells.Item(22, 7).Value = 19480.63
$ws.Cells.Item(35, 6).Value = 91
$ws.Cells.Item(35, 7).Value = 4661.93
$ws.Cells.Item(41, 2).Value = 82525.28
$ws.Cells.Item(43, 6).Value = 162
$ws.Cells.Item(43, 7).Value = 31876.74
$ws.Cells.Item(54, 6).Value = 112
$ws.Cells.Item(54, 7).Value = 2832.48
$ws.Cells.Item(74, 2).Value = 288522.89
$ws.Cells.Item(182, 6).Value = 5
$ws.Cells.Item(182, 7).Value = 328.35
$ws.Cells.Item(184, 2).Value = 30127.13
$ws.Cells.Item(187, 6).Value = 161
$ws.Cells.Item(187, 7).Value = 7167.72
$ws.Cells.Item(192, 2).Value = 49633.25
$ws.Cells.Item(241, 2).Value = 57552
$ws.Cells.Item(241, 5).Value = 136.86
$ws.Cells.Item(241, 6).Value = -5
$ws.Cells.Item(241, 7).Value = -603.45
$ws.Cells.Item(242, 2).Value = 64329
$ws.Cells.Item(242, 5).Value = 128.32
$ws.Cells.Item(242, 6).Value = 1
$ws.Cells.Item(242, 7).Value = 120.69
$ws.Cells.Item(269, 6).Value = 4
$ws.Cells.Item(269, 7).Value = 353.16
$ws.Cells.Item(273, 2).Value = 9638.65
$ws.Cells.Item(278, 6).Value = 252
$ws.Cells.Item(278, 7).Value = 31923.36
$ws.Cells.Item(283, 2).Value = 115664.86
$ws.Cells.Item(316, 6).Value = 14
$ws.Cells.Item(316, 7).Value = 3215.52
$ws.Cells.Item(318, 2).Value = 26259.59
$ws.Cells.Item(349, 6).Value = 137
$ws.Cells.Item(349, 7).Value = 4773.08
$ws.Cells.Item(362, 6).Value = 238
$ws.Cells.Item(362, 7).Value = 11155.06
$ws.Cells.Item(375, 2).Value = 179641.09
$ws.Cells.Item(386, 2).Value = 55373
$ws.Cells.Item(386, 5).Value = 163.62
$ws.Cells.Item(386, 6).Value = -94
$ws.Cells.Item(386, 7).Value = -13562.32
$ws.Cells.Item(387, 2).Value = 63520
$ws.Cells.Item(387, 5).Value = 153.4
$ws.Cells.Item(387, 6).Value = 46
$ws.Cells.Item(387, 7).Value = 6636.88
$ws.Cells.Item(390, 2).Value = 63510
$ws.Cells.Item(390, 5).Value = 50.66
$ws.Cells.Item(390, 6).Value = 86
$ws.Cells.Item(390, 7).Value = 4097.04
$ws.Cells.Item(391, 2).Value = 55356
$ws.Cells.Item(391, 5).Value = 54.04
$ws.Cells.Item(391, 6).Value = -158
$ws.Cells.Item(391, 7).Value = -7527.12
$ws.Cells.Item(400, 2).Value = 63560
$ws.Cells.Item(400, 5).Value = 134.87
$ws.Cells.Item(400, 6).Value = 1
$ws.Cells.Item(400, 7).Value = 126.86
$ws.Cells.Item(401, 2).Value = 60325
$ws.Cells.Item(401, 5).Value = 151.57
$ws.Cells.Item(401, 6).Value = -102
$ws.Cells.Item(401, 7).Value = -12939.72
$ws.Cells.Item(437, 6).Value = 13
$ws.Cells.Item(437, 7).Value = 4439.11
$ws.Cells.Item(450, 6).Value = 138
$ws.Cells.Item(450, 7).Value = 17237.58
$ws.Cells.Item(454, 2).Value = 99408.42
$ws.Cells.Item(468, 6).Value = 111
$ws.Cells.Item(468, 7).Value = 5231.43
$ws.Cells.Item(471, 6).Value = 344
$ws.Cells.Item(471, 7).Value = 57100.56
$ws.Cells.Item(473, 2).Value = 101121.82
$ws.Cells.Item(505, 6).Value = 97
$ws.Cells.Item(505, 7).Value = 3328.07
$ws.Cells.Item(506, 6).Value = 129
$ws.Cells.Item(506, 7).Value = 5229.66
$ws.Cells.Item(522, 2).Value = 207398.57
$ws.Cells.Item(553, 2).Value = 65066
$ws.Cells.Item(553, 5).Value = 13.61
$ws.Cells.Item(553, 6).Value = 90
$ws.Cells.Item(553, 7).Value = 1152.9
$ws.Cells.Item(554, 2).Value = 53263
$ws.Cells.Item(554, 5).Value = 15.29
$ws.Cells.Item(554, 6).Value = -309
$ws.Cells.Item(554, 7).Value = -3958.29
$ws.Cells.Item(572, 2).Value = 53595
$ws.Cells.Item(572, 5).Value = 17.61
$ws.Cells.Item(572, 6).Value = -335
$ws.Cells.Item(572, 7).Value = -4934.55
$ws.Cells.Item(573, 2).Value = 65067
$ws.Cells.Item(573, 5).Value = 15.65
$ws.Cells.Item(573, 6).Value = 126
$ws.Cells.Item(573, 7).Value = 1855.98
$ws.Cells.Item(608, 6).Value = 72
$ws.Cells.Item(608, 7).Value = 19527.84
$ws.Cells.Item(609, 6).Value = 68
$ws.Cells.Item(609, 7).Value = 9881.76
$ws.Cells.Item(615, 2).Value = 150309.01
$ws.Cells.Item(662, 2).Value = 60025
$ws.Cells.Item(662, 5).Value = 37.22
$ws.Cells.Item(662, 6).Value = -98
$ws.Cells.Item(662, 7).Value = -3217.34
$ws.Cells.Item(663, 2).Value = 64833
$ws.Cells.Item(663, 5).Value = 34.9
$ws.Cells.Item(663, 6).Value = 90
$ws.Cells.Item(663, 7).Value = 2954.7
$ws.Cells.Item(679, 6).Value = 246
$ws.Cells.Item(679, 7).Value = 39478.08
$ws.Cells.Item(695, 2).Value = 188405.71
$ws.Cells.Item(706, 6).Value = 63
$ws.Cells.Item(706, 7).Value = 1895.04
$ws.Cells.Item(708, 2).Value = 41669.17
$ws.Cells.Item(731, 6).Value = 8
$ws.Cells.Item(731, 7).Value = 2032.08
$ws.Cells.Item(732, 2).Value = 38422.82
$ws.Cells.Item(735, 6).Value = 328
$ws.Cells.Item(735, 7).Value = 39966.8
$ws.Cells.Item(742, 2).Value = 50674.54
$ws.Cells.Item(744, 6).Value = 39
$ws.Cells.Item(744, 7).Value = 7287.54
$ws.Cells.Item(745, 6).Value = 18
$ws.Cells.Item(745, 7).Value = 3466.08
$ws.Cells.Item(747, 6).Value = 18
$ws.Cells.Item(747, 7).Value = 5876.28
$ws.Cells.Item(750, 6).Value = 32
$ws.Cells.Item(750, 7).Value = 18711.04
$ws.Cells.Item(752, 6).Value = 27
$ws.Cells.Item(752, 7).Value = 15164.55
$ws.Cells.Item(755, 2).Value = 78964.66
$ws.Cells.Item(766, 6).Value = 13
$ws.Cells.Item(766, 7).Value = 8320.780000000001
$ws.Cells.Item(767, 2).Value = 8320.780000000001
$ws.Cells.Item(805, 6).Value = 16
$ws.Cells.Item(805, 7).Value = 2622.24
$ws.Cells.Item(807, 6).Value = 158
$ws.Cells.Item(807, 7).Value = 17191.98
$ws.Cells.Item(811, 6).Value = 296
$ws.Cells.Item(811, 7).Value = 44521.36
$ws.Cells.Item(812, 6).Value = 43
$ws.Cells.Item(812, 7).Value = 6292.19
$ws.Cells.Item(816, 6).Value = 157
$ws.Cells.Item(816, 7).Value = 24283.19
$ws.Cells.Item(817, 6).Value = 81
$ws.Cells.Item(817, 7).Value = 11564.37
$ws.Cells.Item(825, 6).Value = 456
$ws.Cells.Item(825, 7).Value = 35837.04
$ws.Cells.Item(827, 6).Value = 579
$ws.Cells.Item(827, 7).Value = 59584.89
$ws.Cells.Item(832, 6).Value = 98
$ws.Cells.Item(832, 7).Value = 4626.58
$ws.Cells.Item(838, 2).Value = 335170.22
$ws.Cells.Item(890, 6).Value = 249
$ws.Cells.Item(890, 7).Value = 7527.27
$ws.Cells.Item(897, 2).Value = 347353.87
$ws.Cells.Item(929, 6).Value = 135
$ws.Cells.Item(929, 7).Value = 5290.65
$ws.Cells.Item(936, 2).Value = 119467.43
$ws.Cells.Item(939, 6).Value = 134
$ws.Cells.Item(939, 7).Value = 13871.68
$ws.Cells.Item(941, 2).Value = 21205.65
$ws.Cells.Item(942, 2).Value = 5155569.12
$ws.Cells.Item(943, 2).Value = 5155569.12
